$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.091.75'
$ws.Range("E2").Value = '  +5.63%  '
$ws.Range("D3").Value = '1.921.27'
$ws.Range("E3").Value = '  +2.78%  '
$ws.Range("E4").Value = '  -0.72%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.35'
$ws.Range("E5").Value = '  +4.76%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5226'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4100'
$ws.Range("E8").Value = '  +5.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08525'
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.128'
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.82'
$ws.Range("E11").Value = '  +2.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.39'
$ws.Range("E12").Value = '  +9.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.438'
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("D14").Value = '1.905.99'
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.422'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.62'
$ws.Range("E17").Value = '  +4.96%  '
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.42'
$ws.Range("E20").Value = '  +3.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").Value = '30.090.80'
$ws.Range("E23").Value = '  +5.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.34'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D26").Value = '2.139.60'
$ws.Range("E26").Value = '  +2.57%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.12'
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.07'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.454'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.42'
$ws.Range("E30").Value = '  +2.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.083'
$ws.Range("E31").Value = '  +4.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1057'
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.068'
$ws.Range("E33").Value = '  +5.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.629'
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02490'
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06629'
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2220'
$ws.Range("E37").Value = '  +2.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.238'
$ws.Range("E38").Value = '  +4.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.195'
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.895'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6549'
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.64'
$ws.Range("E42").Value = '  +4.95%  '
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6164'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.25'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.771'
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.086'
$ws.Range("E47").Value = '  +4.05%  '
$ws.Range("E48").Value = '  +2.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.68'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.171'
$ws.Range("E50").Value = '  +7.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.78'
$ws.Range("E51").Value = '  +4.10%  '
